$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A75").Value = "2023-12-07 16:19:28"
$ws.Range("B75").Value = 0.003600000000000001
$ws.Range("A76").Value = "2023-12-07 16:19:51"
$ws.Range("B76").Value = 0.0016
